$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would be silently re-typed as a Number by Excel's
# auto-detection (losing a significant trailing zero, e.g. "520.80" -> 520.8);
# force the cell format to Text first so the literal string round-trips.
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"

$ws.Range('D2').Value = '70.904.20'
$ws.Range('E2').Value = '  +6.19%  '
$ws.Range('D3').Value = '3.655.09'
$ws.Range('E3').Value = '  +17.86%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '618.37'
$ws.Range('E5').Value = '  +7.24%  '
$ws.Range('D6').Value = '181.58'
$ws.Range('E6').Value = '  +2.22%  '
$ws.Range('D7').Value = '3.651.95'
$ws.Range('E7').Value = '  +17.82%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.543'
$ws.Range('E9').Value = '  +5.60%  '
$ws.Range('E10').Value = '  +8.49%  '
$ws.Range('E11').Value = '  +5.13%  '
$ws.Range('D12').Value = '0.502'
$ws.Range('E12').Value = '  +7.27%  '
$ws.Range('E13').Value = '  +11.77%  '
$ws.Range('E14').Value = '  +6.08%  '
$ws.Range('D15').Value = '4.259.60'
$ws.Range('E15').Value = '  +17.73%  '
$ws.Range('D16').Value = '70.918.32'
$ws.Range('E16').Value = '  +6.19%  '
$ws.Range('D17').Value = '3.649.18'
$ws.Range('E17').Value = '  +17.52%  '
$ws.Range('E18').Value = '  +1.95%  '
$ws.Range('D19').Value = '7.56'
$ws.Range('E19').Value = '  +7.44%  '
$ws.Range('D20').Value = '520.80'
$ws.Range('E20').Value = '  +8.15%  '
$ws.Range('D21').Value = '16.88'
$ws.Range('E21').Value = '  +1.46%  '
$ws.Range('D22').Value = '9.30'
$ws.Range('E22').Value = '  +18.67%  '
$ws.Range('E23').Value = '  +7.67%  '
$ws.Range('D24').Value = '2.54'
$ws.Range('E24').Value = '  +13.41%  '
$ws.Range('D25').Value = '88.74'
$ws.Range('E25').Value = '  +6.03%  '
$ws.Range('D26').Value = '13.46'
$ws.Range('E26').Value = '  +6.95%  '
$ws.Range('E27').Value = '  +10.10%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').Value = '2.56'
$ws.Range('E29').Value = '  +11.66%  '
$ws.Range('D30').Value = '8.20'
$ws.Range('E30').Value = '  +4.02%  '
$ws.Range('D31').Value = '2.89'
$ws.Range('E31').Value = '  +11.10%  '
$ws.Range('D32').Value = '31.61'
$ws.Range('E32').Value = '  +13.06%  '
$ws.Range('E33').Value = '  +17.43%  '
$ws.Range('E34').Value = '  +4.31%  '
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').Value = '6.13'
$ws.Range('E36').Value = '  +9.41%  '
$ws.Range('B37').Value = 'TheGraph'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D37').Value = '0.352'
$ws.Range('E37').Value = '  +12.87%  '
$ws.Range('B38').Value = 'Mantle'
$ws.Range('C38').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D38').Value = '1.02'
$ws.Range('E38').Value = '  +8.69%  '
$ws.Range('E39').Value = '  +10.06%  '
$ws.Range('D40').Value = '0.131'
$ws.Range('E40').Value = '  +6.31%  '
$ws.Range('D41').Value = '51.35'
$ws.Range('E41').Value = '  +4.66%  '
$ws.Range('D42').Value = '46.07'
$ws.Range('E42').Value = '  -5.04%  '
$ws.Range('B43').Value = 'Cosmos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D43').Value = '8.83'
$ws.Range('E43').Value = '  +5.93%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').Value = '426.26'
$ws.Range('E44').Value = '  +13.57%  '
$ws.Range('D45').Value = '3.116.23'
$ws.Range('E45').Value = '  +11.23%  '
$ws.Range('D46').Value = '2.77'
$ws.Range('E46').Value = '  +3.62%  '
$ws.Range('E47').Value = '  +7.66%  '
$ws.Range('D48').Value = '28.45'
$ws.Range('E48').Value = '  +11.33%  '
$ws.Range('D49').Value = '140.76'
$ws.Range('E49').Value = '  +3.88%  '
$ws.Range('D51').Value = '2.48'
$ws.Range('E51').Value = '  +11.07%  '
